$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "63.118.05"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "2.453.51"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "573.21"
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("D6").Value = "146.18"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").Value = "2.450.44"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  +1.48%  "
$ws.Range("D12").Value = "5.26"
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "0.356"
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "27.01"
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "63.090.08"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "2.446.48"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("E20").Value = "  +5.25%  "
$ws.Range("D21").Value = "328.69"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "4.21"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "2.07"
$ws.Range("E23").Value = "  +13.86%  "
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "65.25"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("D26").Value = "615.35"
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("D27").Value = "8.84"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").Value = "0.0000102"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("D29").Value = "2.580.72"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  +3.95%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").Value = "1.90"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("D35").Value = "5.19"
$ws.Range("E35").Value = "  +6.79%  "
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  +1.29%  "
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "0.379"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "18.86"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").Value = "5.40"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "146.97"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "1.79"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  +5.41%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("D46").Value = "148.88"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("E47").Value = "  +2.72%  "
$ws.Range("D48").Value = "21.14"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("D49").Value = "0.0533"
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("D50").Value = "0.601"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("E51").Value = "  +0.74%  "
